$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 11461.154
$ws.Cells.Item(74, 9).Value = 17713.857
$ws.Cells.Item(74, 10).Value = 4166.3335
$ws.Cells.Item(74, 11).Value = 17713.857
$ws.Cells.Item(74, 12).Value = 4166.3335
$ws.Cells.Item(74, 13).Value = -16777.857
$ws.Cells.Item(74, 14).Value = -6038.3335
$ws.Cells.Item(77, 8).Value = 11461.154
$ws.Cells.Item(77, 9).Value = 17713.857
$ws.Cells.Item(77, 10).Value = 4166.3335
$ws.Cells.Item(77, 11).Value = 88569.285
$ws.Cells.Item(77, 12).Value = 20831.6675
$ws.Cells.Item(77, 13).Value = -83889.285
$ws.Cells.Item(77, 14).Value = -30191.6675
$ws.Cells.Item(100, 8).Value = 29413322
$ws.Cells.Item(100, 9).Value = 1611.6
$ws.Cells.Item(100, 11).Value = 1611.6
$ws.Cells.Item(100, 13).Value = -1070.6
$ws.Cells.Item(113, 8).Value = 22730586
$ws.Cells.Item(113, 9).Value = 45456864
$ws.Cells.Item(113, 10).Value = 4308
$ws.Cells.Item(113, 11).Value = 45456864
$ws.Cells.Item(113, 12).Value = 4308
$ws.Cells.Item(113, 13).Value = -45453610
$ws.Cells.Item(113, 14).Value = -10816
$ws.Cells.Item(141, 8).Value = 1374.8334
$ws.Cells.Item(141, 9).Value = 969.1429000000001
$ws.Cells.Item(141, 10).Value = 4214.6665
$ws.Cells.Item(141, 11).Value = 2907.4287
$ws.Cells.Item(141, 12).Value = 12643.9995
$ws.Cells.Item(141, 13).Value = 2272.5713
$ws.Cells.Item(141, 14).Value = -23003.9995

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3841.25
$ws.Cells.Item(45, 9).Value = 2969.182
$ws.Cells.Item(45, 10).Value = 5759.8
$ws.Cells.Item(45, 11).Value = 2969.182
$ws.Cells.Item(45, 12).Value = 5759.8
$ws.Cells.Item(45, 13).Value = -2592.182
$ws.Cells.Item(45, 14).Value = -6513.8
$ws.Cells.Item(61, 8).Value = 1630.4546
$ws.Cells.Item(61, 9).Value = 1716.1111
$ws.Cells.Item(61, 10).Value = 1245
$ws.Cells.Item(61, 11).Value = 1716.1111
$ws.Cells.Item(61, 12).Value = 1245
$ws.Cells.Item(61, 13).Value = -1504.1111
$ws.Cells.Item(61, 14).Value = -1669
$ws.Cells.Item(97, 8).Value = 2650.238
$ws.Cells.Item(97, 9).Value = 1535.3636
$ws.Cells.Item(97, 11).Value = 1535.3636
$ws.Cells.Item(97, 13).Value = -1039.3636
$ws.Cells.Item(122, 8).Value = 2011.85
$ws.Cells.Item(122, 9).Value = 1959.8422
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 5879.5266
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -3429.5266
$ws.Cells.Item(122, 14).Value = -13900
$ws.Cells.Item(132, 8).Value = 89431.56
$ws.Cells.Item(132, 9).Value = 136647.73
$ws.Cells.Item(132, 10).Value = 2081.65
$ws.Cells.Item(132, 11).Value = 409943.1900000001
$ws.Cells.Item(132, 12).Value = 6244.950000000001
$ws.Cells.Item(132, 13).Value = -407413.1900000001
$ws.Cells.Item(132, 14).Value = -11304.95
$ws.Cells.Item(136, 8).Value = 1630.4546
$ws.Cells.Item(136, 9).Value = 1716.1111
$ws.Cells.Item(136, 10).Value = 1245
$ws.Cells.Item(136, 11).Value = 5148.3333
$ws.Cells.Item(136, 12).Value = 3735
$ws.Cells.Item(136, 13).Value = -2598.3333
$ws.Cells.Item(136, 14).Value = -8835

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 879.5417
$ws.Cells.Item(94, 9).Value = 782.41174
$ws.Cells.Item(94, 11).Value = 782.41174
$ws.Cells.Item(94, 13).Value = -331.41174
$ws.Cells.Item(99, 8).Value = 2489.0908
$ws.Cells.Item(99, 10).Value = 3366.6667
$ws.Cells.Item(99, 12).Value = 3366.6667
$ws.Cells.Item(99, 14).Value = -6362.6667
$ws.Cells.Item(107, 8).Value = 4894.95
$ws.Cells.Item(107, 9).Value = 5575.9165
$ws.Cells.Item(107, 10).Value = 3873.5
$ws.Cells.Item(107, 11).Value = 5575.9165
$ws.Cells.Item(107, 12).Value = 3873.5
$ws.Cells.Item(107, 13).Value = -3655.9165
$ws.Cells.Item(107, 14).Value = -7713.5
$ws.Cells.Item(134, 8).Value = 105024.27
$ws.Cells.Item(134, 9).Value = 126575.414
$ws.Cells.Item(134, 10).Value = 1578.8
$ws.Cells.Item(134, 11).Value = 379726.242
$ws.Cells.Item(134, 12).Value = 4736.4
$ws.Cells.Item(134, 13).Value = -377191.242
$ws.Cells.Item(134, 14).Value = -9806.4

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2098.8
$ws.Cells.Item(58, 9).Value = 2119.3333
$ws.Cells.Item(58, 10).Value = 1914
$ws.Cells.Item(58, 11).Value = 2119.3333
$ws.Cells.Item(58, 12).Value = 1914
$ws.Cells.Item(58, 13).Value = -1916.3333
$ws.Cells.Item(58, 14).Value = -2320
$ws.Cells.Item(99, 8).Value = 1504.0555
$ws.Cells.Item(99, 9).Value = 1421.1305
$ws.Cells.Item(99, 11).Value = 1421.1305
$ws.Cells.Item(99, 13).Value = 76.86950000000002
$ws.Cells.Item(126, 8).Value = 1504.0555
$ws.Cells.Item(126, 9).Value = 1421.1305
$ws.Cells.Item(126, 11).Value = 4263.3915
$ws.Cells.Item(126, 13).Value = -1793.3915
$ws.Cells.Item(132, 8).Value = 4402.1
$ws.Cells.Item(132, 9).Value = 3860.4285
$ws.Cells.Item(132, 10).Value = 5666
$ws.Cells.Item(132, 11).Value = 11581.2855
$ws.Cells.Item(132, 12).Value = 16998
$ws.Cells.Item(132, 13).Value = -9051.2855
$ws.Cells.Item(132, 14).Value = -22058
$ws.Cells.Item(134, 8).Value = 5037.0312
$ws.Cells.Item(134, 9).Value = 5373
$ws.Cells.Item(134, 11).Value = 16119
$ws.Cells.Item(134, 13).Value = -13584
$ws.Cells.Item(136, 8).Value = 2098.8
$ws.Cells.Item(136, 9).Value = 2119.3333
$ws.Cells.Item(136, 10).Value = 1914
$ws.Cells.Item(136, 11).Value = 6357.999899999999
$ws.Cells.Item(136, 12).Value = 5742
$ws.Cells.Item(136, 13).Value = -3807.999899999999
$ws.Cells.Item(136, 14).Value = -10842

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1855.6923
$ws.Cells.Item(102, 9).Value = 1613.7778
$ws.Cells.Item(102, 11).Value = 1613.7778
$ws.Cells.Item(102, 13).Value = 8.22219999999993
$ws.Cells.Item(122, 8).Value = 2145.889
$ws.Cells.Item(122, 9).Value = 2014.5333
$ws.Cells.Item(122, 10).Value = 2802.6667
$ws.Cells.Item(122, 11).Value = 6043.5999
$ws.Cells.Item(122, 12).Value = 8408.000100000001
$ws.Cells.Item(122, 13).Value = -3593.5999
$ws.Cells.Item(122, 14).Value = -13308.0001
$ws.Cells.Item(126, 8).Value = 4053.5293
$ws.Cells.Item(126, 9).Value = 2224.4443
$ws.Cells.Item(126, 10).Value = 6111.25
$ws.Cells.Item(126, 11).Value = 6673.3329
$ws.Cells.Item(126, 12).Value = 18333.75
$ws.Cells.Item(126, 13).Value = -4203.3329
$ws.Cells.Item(126, 14).Value = -23273.75
$ws.Cells.Item(132, 8).Value = 3486.1082
$ws.Cells.Item(132, 9).Value = 3281.923
$ws.Cells.Item(132, 11).Value = 9845.769
$ws.Cells.Item(132, 13).Value = -7315.769

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7144844
$ws.Cells.Item(7, 9).Value = 10001641
$ws.Cells.Item(7, 10).Value = 2850.75
$ws.Cells.Item(7, 11).Value = 10001641
$ws.Cells.Item(7, 12).Value = 2850.75
$ws.Cells.Item(7, 13).Value = -10001529
$ws.Cells.Item(7, 14).Value = -3074.75
$ws.Cells.Item(40, 8).Value = 1627.4375
$ws.Cells.Item(40, 9).Value = 1390.909
$ws.Cells.Item(40, 10).Value = 2147.8
$ws.Cells.Item(40, 11).Value = 1390.909
$ws.Cells.Item(40, 12).Value = 2147.8
$ws.Cells.Item(40, 13).Value = -1254.909
$ws.Cells.Item(40, 14).Value = -2419.8
$ws.Cells.Item(61, 8).Value = 2722.2222
$ws.Cells.Item(61, 9).Value = 1933.3334
$ws.Cells.Item(61, 10).Value = 4300
$ws.Cells.Item(61, 11).Value = 1933.3334
$ws.Cells.Item(61, 12).Value = 4300
$ws.Cells.Item(61, 13).Value = -1731.3334
$ws.Cells.Item(61, 14).Value = -4704
$ws.Cells.Item(93, 8).Value = 2041.7858
$ws.Cells.Item(93, 9).Value = 1871
$ws.Cells.Item(93, 10).Value = 2668
$ws.Cells.Item(93, 11).Value = 1871
$ws.Cells.Item(93, 12).Value = 2668
$ws.Cells.Item(93, 13).Value = -623
$ws.Cells.Item(93, 14).Value = -5164
$ws.Cells.Item(100, 8).Value = 2086.1475
$ws.Cells.Item(100, 9).Value = 2183.4614
$ws.Cells.Item(100, 11).Value = 2183.4614
$ws.Cells.Item(100, 13).Value = -1642.4614
$ws.Cells.Item(113, 8).Value = 2722.2222
$ws.Cells.Item(113, 9).Value = 1933.3334
$ws.Cells.Item(113, 10).Value = 4300
$ws.Cells.Item(113, 11).Value = 1933.3334
$ws.Cells.Item(113, 12).Value = 4300
$ws.Cells.Item(113, 13).Value = 236.6666
$ws.Cells.Item(113, 14).Value = -8640
$ws.Cells.Item(122, 8).Value = 3208.125
$ws.Cells.Item(122, 9).Value = 2366.6667
$ws.Cells.Item(122, 10).Value = 3713
$ws.Cells.Item(122, 11).Value = 7100.000100000001
$ws.Cells.Item(122, 12).Value = 11139
$ws.Cells.Item(122, 13).Value = -4650.000100000001
$ws.Cells.Item(122, 14).Value = -16039
$ws.Cells.Item(126, 8).Value = 7144844
$ws.Cells.Item(126, 9).Value = 10001641
$ws.Cells.Item(126, 10).Value = 2850.75
$ws.Cells.Item(126, 11).Value = 30004923
$ws.Cells.Item(126, 12).Value = 8552.25
$ws.Cells.Item(126, 13).Value = -30002453
$ws.Cells.Item(126, 14).Value = -13492.25
$ws.Cells.Item(132, 8).Value = 2615.7222
$ws.Cells.Item(132, 9).Value = 1956.6666
$ws.Cells.Item(132, 11).Value = 5869.9998
$ws.Cells.Item(132, 13).Value = -3339.9998
$ws.Cells.Item(136, 8).Value = 1502.7368
$ws.Cells.Item(136, 9).Value = 1402.6428
$ws.Cells.Item(136, 10).Value = 1783
$ws.Cells.Item(136, 11).Value = 4207.928400000001
$ws.Cells.Item(136, 12).Value = 5349
$ws.Cells.Item(136, 13).Value = -1657.928400000001
$ws.Cells.Item(136, 14).Value = -10449

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 4522.8667
$ws.Cells.Item(122, 9).Value = 5342.625
$ws.Cells.Item(122, 11).Value = 16027.875
$ws.Cells.Item(122, 13).Value = -13577.875
$ws.Cells.Item(126, 8).Value = 1898.3793
$ws.Cells.Item(126, 9).Value = 1556.5555
$ws.Cells.Item(126, 10).Value = 2457.7273
$ws.Cells.Item(126, 11).Value = 4669.666499999999
$ws.Cells.Item(126, 12).Value = 7373.1819
$ws.Cells.Item(126, 13).Value = -2199.666499999999
$ws.Cells.Item(126, 14).Value = -12313.1819
$ws.Cells.Item(132, 8).Value = 2103.2856
$ws.Cells.Item(132, 9).Value = 1843.7812
$ws.Cells.Item(132, 10).Value = 2933.7
$ws.Cells.Item(132, 11).Value = 5531.3436
$ws.Cells.Item(132, 12).Value = 8801.099999999999
$ws.Cells.Item(132, 13).Value = -3001.3436
$ws.Cells.Item(132, 14).Value = -13861.1
$ws.Cells.Item(136, 8).Value = 1407.9387
$ws.Cells.Item(136, 9).Value = 1282.2162
$ws.Cells.Item(136, 11).Value = 3846.6486
$ws.Cells.Item(136, 13).Value = -1296.6486
